# Update the "dSF" (column F) values for the rows whose underlying source
# data was repulled. Only column F changes; all other columns are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = -1
    5  = -5
    14 = -1
    16 = 0
    19 = -1
    22 = -2
    25 = -3
    27 = -1
    28 = 1
    29 = 1
    35 = -4
    36 = -3
    41 = -1
    45 = 5
    47 = 1
    51 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
